$d = $word.ActiveDocument

# The document is a single paragraph ending in an underlined run whose
# text is "ww". The edit splits that run's text into five separate
# underlined runs: "w", "AAA", "asas", "a", "w" (as if someone placed
# the cursor between the two original w's and kept typing in bursts).
#
# Note on this engine's run handling: any in-place edit that leaves a
# residual text fragment touching an existing run boundary (Delete(),
# Range.Text = "...", or a Find/Replace that shortens text) causes
# adjacent runs that share identical formatting to be coalesced back
# into a single <w:r>, which would wipe out the separate-run structure
# the diff expects. Removing a run *completely* (leaving nothing
# behind) does not trigger that coalescing, and neither does appending
# new text at the very end of the story. So: delete the trailing "ww"
# run outright, then rebuild the desired text as a sequence of freshly
# appended, individually-formatted runs.

$end = $d.Content.End
$tail = $d.Range($end - 3, $end - 1)   # last two chars, excluding the paragraph mark

if ($tail.Text -ne "ww") {
    throw "Expected trailing 'ww' run, found '$($tail.Text)' instead."
}
$tail.Delete()

$pieces = "w", "AAA", "asas", "a", "w"

$endNow = $d.Content.End
$ins = $d.Range($endNow - 1, $endNow - 1)
foreach ($piece in $pieces) {
    $ins.InsertAfter($piece)
    $ins.Font.Underline = 1
    $ins.Collapse(0)
}
